# Updated cryptos list (prices / 1h volume %) per the source commit.
# Values are written with a leading apostrophe so Excel stores them as
# literal text (matching the original inline-string cells) instead of
# silently coercing price strings like "0.120" or "1.10" into numbers
# and dropping the trailing zero / re-formatting them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.551.31"
$ws.Range("E2").Value = "'  -3.48%  "

$ws.Range("D3").Value = "'3.403.16"
$ws.Range("E3").Value = "'  -4.32%  "

$ws.Range("E4").Value = "'  +0.05%  "

$ws.Range("D5").Value = "'579.08"
$ws.Range("E5").Value = "'  -5.04%  "

$ws.Range("D6").Value = "'132.03"
$ws.Range("E6").Value = "'  -9.42%  "

$ws.Range("E7").Value = "'  +0.04%  "

$ws.Range("D8").Value = "'3.401.34"
$ws.Range("E8").Value = "'  -4.35%  "

$ws.Range("E9").Value = "'  -7.51%  "

$ws.Range("D10").Value = "'0.120"
$ws.Range("E10").Value = "'  -10.68%  "

$ws.Range("D11").Value = "'6.93"
$ws.Range("E11").Value = "'  -11.32%  "

$ws.Range("E12").Value = "'  -11.94%  "

$ws.Range("D13").Value = "'3.980.07"
$ws.Range("E13").Value = "'  -4.23%  "

$ws.Range("D14").Value = "'0.0000176"
$ws.Range("E14").Value = "'  -11.61%  "

$ws.Range("E15").Value = "'  -1.97%  "

$ws.Range("D16").Value = "'3.402.46"
$ws.Range("E16").Value = "'  -3.96%  "

$ws.Range("B17").Value = "'Avalanche"
$ws.Range("C17").Value = "'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").Value = "'25.87"
$ws.Range("E17").Value = "'  -11.79%  "

$ws.Range("B18").Value = "'WrappedBTC"
$ws.Range("C18").Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "'64.526.70"
$ws.Range("E18").Value = "'  -3.25%  "

$ws.Range("D19").Value = "'9.43"
$ws.Range("E19").Value = "'  -14.85%  "

$ws.Range("D20").Value = "'5.65"
$ws.Range("E20").Value = "'  -10.39%  "

$ws.Range("D21").Value = "'13.39"
$ws.Range("E21").Value = "'  -9.68%  "

$ws.Range("D22").Value = "'376.88"
$ws.Range("E22").Value = "'  -12.09%  "

$ws.Range("E23").Value = "'  +0.06%  "

$ws.Range("D24").Value = "'0.538"
$ws.Range("E24").Value = "'  -10.91%  "

$ws.Range("D25").Value = "'71.34"
$ws.Range("E25").Value = "'  -8.36%  "

$ws.Range("D26").Value = "'3.535.47"
$ws.Range("E26").Value = "'  -4.30%  "

$ws.Range("D27").Value = "'0.0000103"
$ws.Range("E27").Value = "'  -12.39%  "

$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "'  -0.13%  "

$ws.Range("D29").Value = "'7.08"
$ws.Range("E29").Value = "'  -12.35%  "

$ws.Range("E30").Value = "'  -12.76%  "

$ws.Range("D31").Value = "'7.89"
$ws.Range("E31").Value = "'  -13.64%  "

$ws.Range("D32").Value = "'3.420.00"
$ws.Range("E32").Value = "'  -3.94%  "

$ws.Range("E33").Value = "'  +0.01%  "

$ws.Range("E34").Value = "'  -7.41%  "

$ws.Range("E35").Value = "'  -11.02%  "

$ws.Range("D36").Value = "'170.87"
$ws.Range("E36").Value = "'  -3.92%  "

$ws.Range("D37").Value = "'1.17"
$ws.Range("E37").Value = "'  -14.41%  "

$ws.Range("D38").Value = "'6.59"
$ws.Range("E38").Value = "'  -14.76%  "

$ws.Range("D39").Value = "'1.44"
$ws.Range("E39").Value = "'  -13.43%  "

$ws.Range("D40").Value = "'4.55"
$ws.Range("E40").Value = "'  -14.73%  "

$ws.Range("D41").Value = "'0.0753"
$ws.Range("E41").Value = "'  -9.87%  "

$ws.Range("E42").Value = "'  -8.59%  "

$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "'  +0.01%  "

$ws.Range("D44").Value = "'41.87"
$ws.Range("E44").Value = "'  -8.16%  "

$ws.Range("D45").Value = "'4.23"
$ws.Range("E45").Value = "'  -16.74%  "

$ws.Range("E46").Value = "'  -12.14%  "

$ws.Range("D47").Value = "'1.10"
$ws.Range("E47").Value = "'  -3.39%  "

$ws.Range("D48").Value = "'21.93"
$ws.Range("E48").Value = "'  -6.64%  "

$ws.Range("D49").Value = "'6.45"
$ws.Range("E49").Value = "'  -10.26%  "

$ws.Range("D50").Value = "'2.180.47"
$ws.Range("E50").Value = "'  -6.81%  "

$ws.Range("D51").Value = "'19.78"
$ws.Range("E51").Value = "'  -11.31%  "
